$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 9) into the two new rows
$ws.Range("A9:I9").Copy()
$ws.Range("A10:I10").PasteSpecial(-4122)
$ws.Range("A11:I11").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Add new row 10: LeetCode 2200 - Find All K-Distant Indices in an Array
$ws.Range("A10").Value = 2200
$ws.Range("B10").Value = "Find All K-Distant Indices in an Array"
$ws.Range("C10").Value = "#array #two-pointers"
$ws.Range("D10").Value = "easy"
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 33
$ws.Range("H10").Value = Get-Date -Year 2025 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("I10").Value = Get-Date -Year 2025 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Add new row 11: LeetCode 153 - Find Minimum in Rotated Sorted Array
$ws.Range("A11").Value = 153
$ws.Range("B11").Value = "Find Minimum in Rotated Sorted Array"
$ws.Range("C11").Value = "#array  #binary-search #重点 "
$ws.Range("D11").Value = "medium"
$ws.Range("E11").Value = 6
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = Get-Date -Year 2025 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Range("I11").Value = Get-Date -Year 2025 -Month 6 -Day 24 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Match row heights for new rows (ht=34, same as row 8)
$ws.Range("A10:I11").RowHeight = 34

# Update the view: scroll/selection state
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("I11").Select()

# Update workbook window size
$excel.ActiveWindow.Width = 29100
$excel.ActiveWindow.Height = 14520
